$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the three worker rows (16-18) so that the order is reversed:
# Row 16 becomes DANILSA DE JESUS NIÑO LEAL (CC / 30840158)
# Row 17 stays  KATERINE GOMEZ GONZALEZ    (CC / 1047376062)
# Row 18 becomes JAIRO LUIS CONTRERAS RUIZ (PE / 817674212051994)

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "30840158"
$ws.Range("D16").Value = "DANILSA DE JESUS NIÑO LEAL"
$ws.Range("E16").Value = "1906"

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047376062"
$ws.Range("D17").Value = "KATERINE GOMEZ GONZALEZ"
$ws.Range("E17").Value = "1906"

$ws.Range("B18").Value = "PE"
$ws.Range("C18").Value = "817674212051994"
$ws.Range("D18").Value = "JAIRO LUIS CONTRERAS RUIZ"
$ws.Range("E18").Value = "1906"

$wb.Save()
